# Edit script for PV-Test-02.xlsx
# 1. Change cell C1 text from "Task Name" to "Name"
# 2. Change the active cell selection from C2 to C1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-01")

# Update the header text in C1
$ws.Range("C1").Value = "Name"

# Move the selection to C1 (matches the saved sheetView selection in the diff)
$ws.Activate()
$ws.Range("C1").Select()
